$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:46:20"
$ws1.Range("A3").Value = "Total filas: 33"

# Insert a new data row at position 28 (shifts old rows 28-35 down to 29-36)
$ws1.Rows.Item(28).Insert()

$ws1.Range("A28").Value = "06:46:20"
$ws1.Range("B28").Value = "07:44"
$ws1.Range("C28").Value = "215A_EL PATO"
$ws1.Range("D28").Value = 58
$ws1.Range("E28").Value = "LP1912"

# Append two brand-new rows at the bottom (37, 38)
$ws1.Range("A37").Value = "06:46:20"
$ws1.Range("B37").Value = "08:41"
$ws1.Range("C37").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D37").Value = 115
$ws1.Range("E37").Value = "LP1912"

$ws1.Range("A38").Value = "06:46:20"
$ws1.Range("B38").Value = "08:44"
$ws1.Range("C38").Value = "215C_EL PATO"
$ws1.Range("D38").Value = 118
$ws1.Range("E38").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:46:20"
$ws2.Range("A3").Value = "Total filas: 8"

# Append two brand-new rows at the bottom (12, 13)
$ws2.Range("A12").Value = "06:46:20"
$ws2.Range("B12").Value = "07:44"
$ws2.Range("C12").Value = "215A_EL PATO"
$ws2.Range("D12").Value = 58
$ws2.Range("E12").Value = "LP1912"

$ws2.Range("A13").Value = "06:46:20"
$ws2.Range("B13").Value = "08:44"
$ws2.Range("C13").Value = "215C_EL PATO"
$ws2.Range("D13").Value = 118
$ws2.Range("E13").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:46:20"
$ws3.Range("A3").Value = "Total filas: 3"

# Append one brand-new row at the bottom (8)
$ws3.Range("A8").Value = "06:46:20"
$ws3.Range("B8").Value = "08:36"
$ws3.Range("C8").Value = "215A_LA PLATA"
$ws3.Range("D8").Value = 110
$ws3.Range("E8").Value = "L6173"
